# Word COM-interop script: remove the stale "_GoBack" bookmark left over
# from the previous edit session (Word normally drops this automatically
# the next time the file is saved by the desktop app).
$d = $word.ActiveDocument

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
